$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.387.33'
$ws.Range("E2").Value = '  +3.30%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.619.73'
$ws.Range("E3").Value = '  +1.36%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '568.83'
$ws.Range("E5").Value = '  +5.54%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.34'
$ws.Range("E6").Value = '  +1.55%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.604'
$ws.Range("E8").Value = '  +4.00%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.637.15'
$ws.Range("E9").Value = '  +1.91%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.77'
$ws.Range("E10").Value = '  +0.05%  '

$ws.Range("E11").Value = '  +4.46%  '

$ws.Range("E12").Value = '  +8.68%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.344'
$ws.Range("E13").Value = '  +3.89%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.093.31'
$ws.Range("E14").Value = '  +1.91%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '60.337.96'
$ws.Range("E15").Value = '  +3.36%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '22.04'
$ws.Range("E16").Value = '  +6.43%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000137'
$ws.Range("E17").Value = '  +2.68%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.640.87'
$ws.Range("E18").Value = '  +1.92%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.53'
$ws.Range("E19").Value = '  +1.58%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '339.87'
$ws.Range("E20").Value = '  +1.16%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.39'
$ws.Range("E21").Value = '  +3.41%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.36'
$ws.Range("E22").Value = '  +3.28%  '

$ws.Range("E23").Value = '  +0.06%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.87'
$ws.Range("E24").Value = '  -1.64%  '

$ws.Range("E25").Value = '  +6.88%  '

$ws.Range("E26").Value = '  +4.25%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  +0.64%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.34'
$ws.Range("E28").Value = '  +4.28%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0791'
$ws.Range("E29").Value = '  +8.38%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("E30").Value = '  +0.02%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.70'
$ws.Range("E31").Value = '  +4.00%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.11'
$ws.Range("E32").Value = '  +3.38%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '160.20'
$ws.Range("E33").Value = '  +3.17%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '19.09'
$ws.Range("E34").Value = '  +1.30%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.08'
$ws.Range("E35").Value = '  +4.78%  '

$ws.Range("E36").Value = '  +4.55%  '

$ws.Range("B37").Value = 'Fetch.AI'
$ws.Range("C37").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.885'
$ws.Range("E37").Value = '  +8.05%  '

$ws.Range("B38").Value = 'SuiNetwork'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.878'
$ws.Range("E38").Value = '  +6.18%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '37.51'
$ws.Range("E39").Value = '  +1.48%  '

$ws.Range("E40").Value = '  +6.80%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '297.41'
$ws.Range("E41").Value = '  +5.30%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.64'
$ws.Range("E42").Value = '  +1.48%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.997'
$ws.Range("E43").Value = '  -0.23%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0985'
$ws.Range("E44").Value = '  +4.60%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.600'
$ws.Range("E45").Value = '  +1.71%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0542'
$ws.Range("E46").Value = '  +1.42%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '19.32'
$ws.Range("E47").Value = '  +4.56%  '

$ws.Range("B48").Value = 'WhiteBITCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '10.70'
$ws.Range("E48").Value = '  +0.69%  '

$ws.Range("B49").Value = 'Aave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '127.20'
$ws.Range("E49").Value = '  +16.00%  '

$ws.Range("E50").Value = '  +3.76%  '

$ws.Range("B51").Value = 'Maker'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.960.49'
$ws.Range("E51").Value = '  +2.44%  '
